$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- 1. Collapse the three CORE COMPETENCIES paragraphs into one summary paragraph ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Data Visualization & Design $bullet Geospatial Analysis & Mapping $bullet Technical Visualization"

# Remove the old second and third Core Competencies paragraphs (now items 7 and 8)
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$oldRange = $d.Range($p7.Range.Start, $p8.Range.End)
$oldRange.Delete()

# --- 2. Append a new "TECHNICAL SKILLS" section at the end of the document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading2"

$headingPara.Range.InsertParagraphAfter()
$body1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$body1.Style = "Normal"
$body1.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design; Web Visualization; Presentation Design; Data Storytelling"

$body1.Range.InsertParagraphAfter()
$body2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$body2.Style = "Normal"
$body2.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing; Census Data Integration; Custom Tile Servers; Spatial Clustering"

$body2.Range.InsertParagraphAfter()
$body3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$body3.Style = "Normal"
$body3.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Cloud Platforms; Web Technologies; Statistical Computing; Version Control; DevOps"
